$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '47.182.00'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.90%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.494.43'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.84%  '

# Row 4
$ws.Range('E4').Value = '  -0.21%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.86%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '105.87'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.30%  '

# Row 7
$ws.Range('E7').Value = '  +1.02%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.20%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.541'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.80%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.57'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.26%  '

# Row 11
$ws.Range('E11').Value = '  +0.90%  '

# Row 12
$ws.Range('E12').Value = '  +0.26%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.33'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.37%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.18'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.17%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.882.87'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.69%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.521.54'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.57%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.844'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.75%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '47.078.15'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.91%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.71'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.28%  '

# Row 20
$ws.Range('E20').Value = '  +2.41%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0935'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.53%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '70.88'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.33%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '251.47'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.75%  '

# Row 24
$ws.Range('E24').Value = '  +2.93%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.54'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.29%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.13'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.28%  '

# Row 27
$ws.Range('E27').Value = '  +0.05%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.15'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.91%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.20'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.38%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.35'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.49%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.134'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.88%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.65'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.27%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.80'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.16%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.38'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.32%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0782'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.96%  '

# Row 36
$ws.Range('E36').Value = '  -0.03%  '

# Row 37
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.63'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.55%  '

# Row 38
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.93'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.96%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.98'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.47%  '

# Row 40
$ws.Range('E40').Value = '  +1.18%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '121.48'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.88%  '

# Row 42
$ws.Range('E42').Value = '  +0.80%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.66'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.26%  '

# Row 44
$ws.Range('E44').Value = '  +0.84%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.952.92'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.54%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.97'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.01%  '

# Row 47
$ws.Range('E47').Value = '  -0.23%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.17'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.03%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.79'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.28%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.38'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +13.11%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.02'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.41%  '
